$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "NEBNextPoly(A)E7490" -> "NEBNextPoly(A)E7490L" (shared string used by G2:G37)
$ws.Range("G2:G37").Value = "NEBNextPoly(A)E7490L"

# 2. Widen column G (polyAIsolationProtocol) so it gets its own <col> span,
#    splitting the former single 1-1025 run into 1-6 / 7-7 / 8-1025.
$ws.Columns.Item(7).ColumnWidth = 27.1

# 3. Move the selection from I2:I37 to G2:G37 (view also scrolls back to A1).
$ws.Range("G2:G37").Select()
